# Generate Report for Handback
#
# For the "045f81b5-f2c0-42ba-bca5-0427652a4527" handback row (row 6) on both
# the zh-cn and de-de status sheets, the handback file turned out to be based
# on a stale commit of the source .md. Populate the "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# columns for that row, add the corresponding hyperlink on the target-file
# cell, and widen the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8aaa1151693811f3120778deef7fac4936916bc5/e2e/045f81b5-f2c0-42ba-bca5-0427652a4527.md"
$latestTargetDisplay = "045f81b5-f2c0-42ba-bca5-0427652a4527.md"

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32f13dfa247e4570e6238916f475e576b2bdd3b6/e2e/045f81b5-f2c0-42ba-bca5-0427652a4527.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8aaa1151693811f3120778deef7fac4936916bc5/e2e/045f81b5-f2c0-42ba-bca5-0427652a4527.md."

$errorDetailDate = "2016-08-24 10:44:32"

function Update-HandbackRow($sheetName, $handbackFile, $handbackDatetimeValue) {
    $ws = $wb.Worksheets.Item($sheetName)

    # I6 = Latest Target File -> hyperlink to the latest commit of the source .md
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestTargetUrl, "", "", $latestTargetDisplay) | Out-Null

    # J6 = Latest Handback File
    $ws.Range("J6").Value = $handbackFile

    # K6 = Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDatetimeValue

    # P6 = Error Detail
    $ws.Range("P6").Value = $errorDetailDate

    # Widen the Error Detail column (P / 16th column) so the message is readable
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-HandbackRow "zh-cn" "045f81b5-f2c0-42ba-bca5-0427652a4527.bd00f40fcf54e02cf6b1a8bba5868754e2b4438b.zh-cn.xlf" $errorMessage
Update-HandbackRow "de-de" "045f81b5-f2c0-42ba-bca5-0427652a4527.bd00f40fcf54e02cf6b1a8bba5868754e2b4438b.de-de.xlf" "2016-08-24 10:44:55"
